# Update benchmark: 2026-01-04 06:40:59 UTC
# Clears several stale duplicate cell values and updates two computed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (HESAPTAN EFT - Sube)
$ws.Range("F3").Value = ""
$ws.Range("I3").Value = ""

# Row 4 (HESAPTAN EFT - ATM)
$ws.Range("F4").Value = ""
$ws.Range("I4").Value = ""

# Row 5 (HESAPTAN EFT - Mobil)
$ws.Range("F5").Value = ""
$ws.Range("I5").Value = ""

# Row 6 (DUZENLI EFT)
$ws.Range("I6").Value = ""

# Row 8 (HESAPTAN HAVALE - Sube)
$ws.Range("F8").Value = ""
$ws.Range("I8").Value = ""

# Row 9 (HESAPTAN HAVALE - ATM)
$ws.Range("F9").Value = ""
$ws.Range("I9").Value = ""

# Row 10 (HESAPTAN HAVALE - Mobil)
$ws.Range("F10").Value = ""
$ws.Range("I10").Value = ""

# Row 11 (DUZENLI HAVALE)
$ws.Range("I11").Value = ""

# Row 13 (GELEN SWIFT)
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("F13").Value = ""
$ws.Range("I13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14 (GIDEN SWIFT - Mobil)
$ws.Range("F14").Value = ""
